$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Set property's "Private" column (D) values to TRUE for the data rows (2-6)
$ws.Range("D2:D6").Value = $true

# Update selection to match the newly-edited column
$ws.Range("D2:D6").Select()
